# Daily attendance processing - 2026-01-06 03:30:49
# Normalize the "Recorded By" (column G) entries: for multi-author cells
# whose first listed recorder is "System" or "backup@backdoor.com",
# rotate that leading entry to the end of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ',\s*'
    if ($parts.Length -gt 1 -and ($parts[0] -eq 'System' -or $parts[0] -eq 'backup@backdoor.com')) {
        $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ', '
        $cell.Value = $rotated
    }
}
